# Apply the edits described by the commit:
# "WIP: Improve gen_top logic to preserve original module structure and add parity ports/instance"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAFETY.PARITY")

# --- Update BIT WIDTH (P) / PARITY SOURCE BIT WIDTH (Q) / MD5 & Script Version (U) ---

$newMd5 = "MD5: 0361ade3b14b1220359a6cc5e306a61e | Script: v3.0.0"

# Row 3
$ws.Range("P3").Value = 35
$ws.Range("Q3").Value = 8
$ws.Range("U3").Value = $newMd5

# Row 4
$ws.Range("P4").Value = 35
$ws.Range("Q4").Value = 8
$ws.Range("U4").Value = $newMd5

# Row 5 (P unchanged, 256)
$ws.Range("Q5").Value = 8
$ws.Range("U5").Value = $newMd5

# Row 6 (P unchanged, 256)
$ws.Range("Q6").Value = 8
$ws.Range("U6").Value = $newMd5

# Row 7
$ws.Range("P7").Value = 34
$ws.Range("Q7").Value = 8
$ws.Range("U7").Value = $newMd5

# Row 8
$ws.Range("P8").Value = 34
$ws.Range("Q8").Value = 8
$ws.Range("U8").Value = $newMd5

# Row 9 (P unchanged, 256)
$ws.Range("Q9").Value = 8
$ws.Range("U9").Value = $newMd5

# Row 10 (P unchanged, 256)
$ws.Range("Q10").Value = 8
$ws.Range("U10").Value = $newMd5

# --- Update the active selection on the sheet (was U2, now Q13) ---
$ws.Range("Q13").Select()

# --- Reposition the workbook window (best effort; mirrors xWindow/yWindow change) ---
$win = $wb.Windows.Item(1)
$win.Left = -28920
$win.Top = -210

# --- Add page setup info (orientation=portrait, paperSize=9/A4) ---
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9

$wb.Save()
